$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-22, columns B-F (A stays as index 0..20, unchanged)
$data = @{
    2  = @{ B = "NSE:ADANIENT";    C = "NSE:DATAPATTNS"; D = $null; E = "NSE:BSOFT"; F = "NSE:MARUTI" }
    3  = @{ B = "NSE:ADROITINFO";  C = "NSE:FDC";         D = $null; E = "NSE:CANBK"; F = $null }
    4  = @{ B = "NSE:ALPA";        C = "NSE:FINPIPE";     D = $null; E = "NSE:MPHASIS"; F = $null }
    5  = @{ B = "NSE:AYMSYNTEX";   C = "NSE:GRASIM";      D = $null; E = "NSE:OFSS"; F = $null }
    6  = @{ B = "NSE:BROOKS";      C = "NSE:ICIL";        D = $null; E = $null; F = $null }
    7  = @{ B = "NSE:CYBERMEDIA";  C = "NSE:IDEAFORGE";   D = $null; E = $null; F = $null }
    8  = @{ B = "NSE:DEEPINDS";    C = "NSE:IKIO";        D = $null; E = $null; F = $null }
    9  = @{ B = "NSE:DOLATALGO";   C = "NSE:INCREDIBLE";  D = $null; E = $null; F = $null }
    10 = @{ B = "NSE:HDFCNIFTY";   C = "NSE:INDHOTEL";    D = $null; E = $null; F = $null }
    11 = @{ B = "NSE:HINDOILEXP";  C = "NSE:IRCON";       D = $null; E = $null; F = $null }
    12 = @{ B = "NSE:JSWSTEEL";    C = "NSE:JAYBARMARU";  D = $null; E = $null; F = $null }
    13 = @{ B = "NSE:JUBLPHARMA";  C = "NSE:KAYNES";      D = $null; E = $null; F = $null }
    14 = @{ B = "NSE:MANORG";      C = "NSE:KHAITANLTD";  D = $null; E = $null; F = $null }
    15 = @{ B = "NSE:MARUTI";      C = "NSE:KRISHANA";    D = $null; E = $null; F = $null }
    16 = @{ B = "NSE:MASTEK";      C = "NSE:MAHLOG";      D = $null; E = $null; F = $null }
    17 = @{ B = "NSE:MATRIMONY";   C = "NSE:MAHSEAMLES";  D = $null; E = $null; F = $null }
    18 = @{ B = "NSE:MAYURUNIQ";   C = "NSE:MAITHANALL";  D = $null; E = $null; F = $null }
    19 = @{ B = "NSE:NIFTYETF";    C = "NSE:MBAPL";       D = $null; E = $null; F = $null }
    20 = @{ B = "NSE:RADHIKAJWE";  C = "NSE:MHRIL";       D = $null; E = $null; F = $null }
    21 = @{ B = "NSE:RELINFRA";    C = "NSE:RHL";         D = $null; E = $null; F = $null }
    22 = @{ B = "NSE:RPOWER";      C = "NSE:ROHLTD";      D = $null; E = $null; F = $null }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    foreach ($col in $cols.Keys) {
        $value = $cols[$col]
        $addr = "$col$row"
        if ($null -eq $value) {
            $ws.Range($addr).Value = ""
        } else {
            $ws.Range($addr).Value = $value
        }
    }
}

# Remove rows 23 to 33 (previously had data, now should be gone entirely)
$ws.Range("A23:F33").EntireRow.Delete()
